# AFDP-4502: Append domain to user/group ids in config files
#
# The "owning group" default-participant rules on Sheet1 hard-coded bare
# LDAP group names. Append the "@ARMEDIA.COM" domain suffix to those two
# values (used by the Case File / Complaint "Default group" rule and the
# Organization / Person "Default group" rule).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Case File / Complaint - Default group action (G23, G24)
$ws.Range("G23").Value = "owning group, ACM_INVESTIGATOR_DEV@ARMEDIA.COM"
$ws.Range("G24").Value = "owning group, ACM_INVESTIGATOR_DEV@ARMEDIA.COM"

# Organization / Person - Default group action (G32, G35)
$ws.Range("G32").Value = "owning group, ACM_ADMINISTRATOR_DEV@ARMEDIA.COM"
$ws.Range("G35").Value = "owning group, ACM_ADMINISTRATOR_DEV@ARMEDIA.COM"

# Widen column G so the longer values with the appended domain still fit,
# and move the active selection to where the author ended up (G35).
$ws.Columns.Item(7).ColumnWidth = 51.7109375

$ws.Range("G35").Select()
